# The two species records (rows 3 and 4) got swapped: the record that
# used to be on row 3 ("Läderdoftande fingersvamp" / Ramaria safraniolens)
# moved to row 4, and the record that used to be on row 4 ("Droppklibbskivling"
# / Limacella guttata) moved to row 3.
#
# Only touch the columns that actually hold per-record data that differs
# between the two rows, leaving everything else (including genuinely blank
# cells) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "B", "D", "E", "F", "G", "H", "P", "Q", "R")

foreach ($col in $columns) {
    $cell3 = $ws.Range($col + "3")
    $cell4 = $ws.Range($col + "4")

    $v3 = $cell3.Value()
    $v4 = $cell4.Value()

    $cell3.Value = $v4
    $cell4.Value = $v3
}
